$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 12:55"

# Update country data: refreshed totals cause several countries to swap
# positions in the (descending-by-total) list, and a few new highs move
# rows up. All literal cell values below reflect the refreshed dataset.

$ws.Range("F6").Value = 5231
$ws.Range("D17").Value = 636
$ws.Range("E17").Value = 8359
$ws.Range("F17").Value = 193
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 108
$ws.Range("F25").Value = 52
$ws.Range("B55").Value = 715
$ws.Range("C55").Value = 2
$ws.Range("E55").Value = 654
$ws.Range("B69").Value = 480
$ws.Range("C69").Value = 5
$ws.Range("E69").Value = 463
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 11
$ws.Range("A75").Value = 'Eslovaquia'
$ws.Range("B75").Value = 336
$ws.Range("C75").Value = 22
$ws.Range("D75").Value = 7
$ws.Range("E75").Value = 329
$ws.Range("F75").Value = 1
$ws.Range("H75").Value = 0
$ws.Range("A76").Value = 'Principado de Andorra'
$ws.Range("B76").Value = 334
$ws.Range("D76").Value = 6
$ws.Range("E76").Value = 322
$ws.Range("F76").Value = 10
$ws.Range("H76").Value = 6
$ws.Range("E85").Value = 237
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 4
$ws.Range("A97").Value = 'Malta'
$ws.Range("B97").Value = 156
$ws.Range("C97").Value = 5
$ws.Range("E97").Value = 154
$ws.Range("F97").Value = 4
$ws.Range("H97").Value = 0
$ws.Range("A98").Value = 'Ghana'
$ws.Range("B98").Value = 152
$ws.Range("E98").Value = 145
$ws.Range("F98").Value = 1
$ws.Range("H98").Value = 5
$ws.Range("A104").Value = 'Sri Lanka'
$ws.Range("B104").Value = 122
$ws.Range("C104").Value = 5
$ws.Range("D104").Value = 15
$ws.Range("E104").Value = 106
$ws.Range("F104").Value = 5
$ws.Range("H104").Value = 1
$ws.Range("A105").Value = 'Afganistan'
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 2
$ws.Range("E105").Value = 114
$ws.Range("F105").Value = 0
$ws.Range("H105").Value = 4
$ws.Range("A118").Value = 'Mayotte'
$ws.Range("B118").Value = 82
$ws.Range("C118").Value = 19
$ws.Range("D118").Value = 10
$ws.Range("E118").Value = 72
$ws.Range("F118").Value = 3
$ws.Range("H118").Value = 0
$ws.Range("A119").Value = 'Consejo Danes para los Refugiados'
$ws.Range("B119").Value = 81
$ws.Range("D119").Value = 2
$ws.Range("E119").Value = 71
$ws.Range("H119").Value = 8
$ws.Range("A120").Value = 'Trinidad yTobago'
$ws.Range("B120").Value = 78
$ws.Range("D120").Value = 1
$ws.Range("E120").Value = 74
$ws.Range("H120").Value = 3
$ws.Range("A121").Value = 'Ruanda'
$ws.Range("B121").Value = 70
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 70
$ws.Range("A122").Value = 'Gibraltar'
$ws.Range("B122").Value = 65
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 14
$ws.Range("E122").Value = 51
$ws.Range("F122").Value = 0
$ws.Range("H122").Value = 0
$ws.Range("A123").Value = 'Paraguay'
$ws.Range("B123").Value = 64
$ws.Range("C123").Value = 5
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = 60
$ws.Range("F123").Value = 3
$ws.Range("H123").Value = 3
$ws.Range("A134").Value = 'Polinesia Francesa'
$ws.Range("B134").Value = 35
$ws.Range("C134").Value = 5
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 35
$ws.Range("F134").Value = 2
$ws.Range("H134").Value = 0
$ws.Range("A135").Value = 'Guatemala'
$ws.Range("B135").Value = 34
$ws.Range("D135").Value = 10
$ws.Range("E135").Value = 23
$ws.Range("F135").Value = 1
$ws.Range("H135").Value = 1
$ws.Range("A136").Value = 'Uganda'
$ws.Range("A137").Value = 'Barbados'
$ws.Range("B137").Value = 33
$ws.Range("E137").Value = 33
$ws.Range("H137").Value = 0
$ws.Range("A138").Value = 'Guam'
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 31
$ws.Range("A139").Value = 'Jamaica'
$ws.Range("B139").Value = 32
$ws.Range("D139").Value = 2
$ws.Range("E139").Value = 29
$ws.Range("H139").Value = 1
$ws.Range("A145").Value = 'Niger'
$ws.Range("B145").Value = 22
$ws.Range("C145").Value = 4
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 3
$ws.Range("A146").Value = 'Congo'
$ws.Range("B146").Value = 19
$ws.Range("E146").Value = 19
$ws.Range("A147").Value = 'Republica de Yibuti'
$ws.Range("E147").Value = 18
$ws.Range("H147").Value = 0
$ws.Range("A148").Value = 'Mali'
$ws.Range("E148").Value = 17
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1
$ws.Range("A154").Value = 'Tanzania'
$ws.Range("A155").Value = 'Bahamas'
$ws.Range("A156").Value = 'Guinea Ecuatorial'
$ws.Range("A157").Value = 'Eritrea'
$ws.Range("A160").Value = 'Dominica'
$ws.Range("A161").Value = 'San Martin (Parte Francesa)'
$ws.Range("A165").Value = 'Granada'
$ws.Range("A166").Value = 'Suazilandia'
$ws.Range("A167").Value = 'Siria'
$ws.Range("D167").Value = 0
$ws.Range("H167").Value = 1
$ws.Range("A168").Value = 'Santa Lucia'
$ws.Range("D168").Value = 1
$ws.Range("H168").Value = 0

$wb.Application.Calculate()
